$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("D12").Value = 8
$ws.Range("D14").Value = 10
$ws.Range("F14").Value = 3
$ws.Range("D19").Value = 9

# Update the selected cell (active cell) to K5
$ws.Range("K5").Select()
